$wb = $excel.ActiveWorkbook

# --- Update test data: rename device "PCH800" to "PCH800 5.0A" on the
# "Add Devices" sheet (cell B2). This is the actual data change behind
# the commit "Updated test data for BVT and TC_66".
$wsDevices = $wb.Worksheets.Item("Add Devices")
$wsDevices.Range("B2").Value = "PCH800 5.0A"

# --- Selection / active-sheet bookkeeping: the workbook now opens with
# "Add Devices" as the active tab, with B2 selected (the cell that was
# just edited). "Add Panels" is no longer the selected tab.
$wsDevices.Range("B2").Select() | Out-Null
$wsDevices.Activate() | Out-Null
